# Rewrite the "Questions" worksheet so the JSON-ish questions payload
# (previously duplicated: a stray 0 in A1 + the real text in A2) lives
# solely in A1, pretty-printed, with A1 restored to default (unstyled)
# formatting and the old helper row removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$questionsText = @'
questions = [
    {
        "title": "In your organization, there are 50 ESXi hosts managed by the vCenter Server. You want to have centralized, automated patch and version management.  Which component will allow you to achieve this goal?",
        "ques_type": 2,
        "options": [
            "VMware vSphere Update Manager",
            "VMware vSphere AutoDeploy",
            "VMware vRealize Orchestrator",
            "VMware vRealize Operations Manager"
        ],
        "score": "VMware vSphere Update Manager"
    },
    {
        "title": "In your work as a VMware administrator, you notice the below alarm on one of your virtual machines (VMs). The VM is working properly, and it is accessible to the users. To resolve this issue, you consolidated the virtual machine, but it failed to consolidate.  What is the best approach to resolve this issue?",
        "ques_type": 2,
        "options": [
            "Ignore the alarm since the VM is working properly.",
            "Create new Snapshots to save the current state of the VM.",
            "Check for VM locks to the virtual machine disk (VMDK) file and resolve it first.",
            "Power down the VM as soon as possible."
        ],
        "score": "Check for VM locks to the virtual machine disk (VMDK) file and resolve it first."
    },
    {
        "title": "During an ESXi host disconnection issue, you want to obtain troubleshooting information from the log files.   Which log files do you check in an ESXi host?",
        "ques_type": 2,
        "options": [
            "/vmfs/volumes/&ltdatastore_name&gt/&ltVM_name&gt/vmware.log",
            "/usr/lib/vmware/hosts.txt",
            "/etc/log/vpxa.log",
            "/var/log/vmware/hostd.log"
        ],
        "score": "/var/log/vmware/hostd.log"
    },
    {
        "title": "To provide network connectivity to virtual machines, what are the available network adapter (vNIC) types in VMware? Select all that apply",
        "ques_type": 15,
        "options": [
            "PCI NIC",
            "E5000",
            "E1000",
            "BNC or RJ45",
            "VMXNET 3"
        ],
        "score": [
            "E1000",
            "VMXNET 3"
        ]
    }
]
'@

# A1 previously held a bold/bordered/centered placeholder value (0).
# Strip that formatting back to the workbook default before writing
# the real text into it.
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $questionsText

# The text used to live in A2 (plain style) - that row is no longer
# needed now that A1 carries the (reformatted) text itself.
$ws.Range("A2").EntireRow.Delete()

# Let row 1 settle back to its natural (default) height now that its
# content changed, rather than leaving it pinned to a stale size.
$ws.Rows.Item(1).AutoFit()
